# Auto-generated script to refresh market-price data cells (columns H-N)
# across the Leve profit sheets, per the scheduled runner update.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 3762.6365
$ws.Range("I98").Value = 3821.8333
$ws.Range("J98").Value = 3496.25
$ws.Range("K98").Value = 3821.8333
$ws.Range("L98").Value = 3496.25
$ws.Range("M98").Value = -2323.8333
$ws.Range("N98").Value = -6492.25
$ws.Range("H107").Value = 1839.2727
$ws.Range("I107").Value = 1604.625
$ws.Range("K107").Value = 1604.625
$ws.Range("M107").Value = 315.375
$ws.Range("H113").Value = 6574.7646
$ws.Range("I113").Value = 8889.454
$ws.Range("J113").Value = 2331.1667
$ws.Range("K113").Value = 8889.454
$ws.Range("L113").Value = 2331.1667
$ws.Range("M113").Value = -5635.454
$ws.Range("N113").Value = -8839.1667
$ws.Range("H122").Value = 3762.6365
$ws.Range("I122").Value = 3821.8333
$ws.Range("J122").Value = 3496.25
$ws.Range("K122").Value = 11465.4999
$ws.Range("L122").Value = 10488.75
$ws.Range("M122").Value = -9015.499899999999
$ws.Range("N122").Value = -15388.75
$ws.Range("H137").Value = 15364.866
$ws.Range("I137").Value = 1406.5454
$ws.Range("K137").Value = 4219.6362
$ws.Range("M137").Value = -1669.6362
$ws.Range("H141").Value = 5226.3335
$ws.Range("I141").Value = 4711
$ws.Range("K141").Value = 14133
$ws.Range("M141").Value = -8953

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H39").Value = 5243.8887
$ws.Range("I39").Value = 4524.375
$ws.Range("J39").Value = 11000
$ws.Range("K39").Value = 4524.375
$ws.Range("L39").Value = 11000
$ws.Range("M39").Value = -4004.375
$ws.Range("N39").Value = -12040
$ws.Range("H132").Value = 2008109.9
$ws.Range("I132").Value = 2080.0232
$ws.Range("K132").Value = 6240.069600000001
$ws.Range("M132").Value = -3710.069600000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("M23").ClearContents()
$ws.Range("N23").ClearContents()
$ws.Range("H134").Value = 49591.605
$ws.Range("J134").Value = 21414.438
$ws.Range("L134").Value = 64243.314
$ws.Range("N134").Value = -69313.314

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 4829.6665
$ws.Range("I2").Value = 4829.6665
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 4829.6665
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -4716.6665
$ws.Range("N2").ClearContents()
$ws.Range("H22").Value = 2012.7273
$ws.Range("I22").Value = 1400
$ws.Range("K22").Value = 1400
$ws.Range("M22").Value = -1050
$ws.Range("H26").Value = 19000
$ws.Range("J26").Value = 19000
$ws.Range("L26").Value = 19000
$ws.Range("N26").Value = -19574
$ws.Range("H31").Value = 13074.154
$ws.Range("J31").Value = 36098.777
$ws.Range("L31").Value = 36098.777
$ws.Range("N31").Value = -36688.777
$ws.Range("H34").Value = 13074.154
$ws.Range("J34").Value = 36098.777
$ws.Range("L34").Value = 36098.777
$ws.Range("N34").Value = -36502.777

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 3537529
$ws.Range("J4").Value = 3000175
$ws.Range("L4").Value = 9000525
$ws.Range("N4").Value = -9000749
$ws.Range("H7").Value = 125.92308
$ws.Range("I7").Value = 125.92308
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 377.76924
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -265.76924
$ws.Range("N7").ClearContents()
$ws.Range("H33").Value = 400.42856
$ws.Range("J33").Value = 450.66666
$ws.Range("L33").Value = 2703.99996
$ws.Range("N33").Value = -3269.99996
$ws.Range("H47").Value = 413
$ws.Range("I47").Value = 503
$ws.Range("J47").Value = 323
$ws.Range("K47").Value = 1509
$ws.Range("L47").Value = 969
$ws.Range("M47").Value = -1078
$ws.Range("N47").Value = -1831
$ws.Range("H131").Value = 1360.56
$ws.Range("I131").Value = 596.3
$ws.Range("J131").Value = 1445.4778
$ws.Range("K131").Value = 1788.9
$ws.Range("L131").Value = 4336.4334
$ws.Range("M131").Value = 3251.1
$ws.Range("N131").Value = -14416.4334
$ws.Range("H133").Value = 4293.143
$ws.Range("I133").Value = 3203.353
$ws.Range("K133").Value = 9610.059000000001
$ws.Range("M133").Value = -4550.059000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 140.04347
$ws.Range("I2").Value = 101.82353
$ws.Range("J2").Value = 248.33333
$ws.Range("K2").Value = 101.82353
$ws.Range("L2").Value = 248.33333
$ws.Range("M2").Value = 11.17646999999999
$ws.Range("N2").Value = -474.33333
$ws.Range("H14").Value = 1116640.8
$ws.Range("I14").Value = 2003172.4
$ws.Range("K14").Value = 2003172.4
$ws.Range("M14").Value = -2003004.4
$ws.Range("H80").Value = 7753.8335
$ws.Range("I80").Value = 6475.909
$ws.Range("J80").Value = 9762
$ws.Range("K80").Value = 6475.909
$ws.Range("L80").Value = 9762
$ws.Range("M80").Value = -5477.909
$ws.Range("N80").Value = -11758
$ws.Range("H83").Value = 7753.8335
$ws.Range("I83").Value = 6475.909
$ws.Range("J83").Value = 9762
$ws.Range("K83").Value = 32379.545
$ws.Range("L83").Value = 48810
$ws.Range("M83").Value = -27387.545
$ws.Range("N83").Value = -58794
$ws.Range("H123").Value = 62000
$ws.Range("J123").Value = 59000
$ws.Range("L123").Value = 59000
$ws.Range("N123").Value = -63900
$ws.Range("H132").Value = 6145.1875
$ws.Range("J132").Value = 13588.111
$ws.Range("L132").Value = 40764.333
$ws.Range("N132").Value = -45824.333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 111114650
$ws.Range("I16").Value = 200003680
$ws.Range("K16").Value = 200003680
$ws.Range("M16").Value = -200003510
$ws.Range("H46").Value = 2369.0435
$ws.Range("I46").Value = 500
$ws.Range("J46").Value = 2454
$ws.Range("K46").Value = 500
$ws.Range("L46").Value = 2454
$ws.Range("M46").Value = -312
$ws.Range("N46").Value = -2830
$ws.Range("H55").Value = 1846.7727
$ws.Range("I55").Value = 1409.5385
$ws.Range("K55").Value = 1409.5385
$ws.Range("M55").Value = -1236.5385
$ws.Range("H61").Value = 1642.0476
$ws.Range("J61").Value = 1601.5
$ws.Range("L61").Value = 1601.5
$ws.Range("N61").Value = -2005.5
$ws.Range("H82").Value = 2675.8572
$ws.Range("I82").Value = 3299.75
$ws.Range("J82").Value = 1844
$ws.Range("K82").Value = 3299.75
$ws.Range("L82").Value = 1844
$ws.Range("M82").Value = -2938.75
$ws.Range("N82").Value = -2566
$ws.Range("H85").Value = 2675.8572
$ws.Range("I85").Value = 3299.75
$ws.Range("J85").Value = 1844
$ws.Range("K85").Value = 3299.75
$ws.Range("L85").Value = 1844
$ws.Range("M85").Value = -2051.75
$ws.Range("N85").Value = -4340
$ws.Range("H113").Value = 1642.0476
$ws.Range("J113").Value = 1601.5
$ws.Range("L113").Value = 1601.5
$ws.Range("N113").Value = -5941.5
$ws.Range("H136").Value = 19966
$ws.Range("J136").Value = 16638.691
$ws.Range("L136").Value = 49916.073
$ws.Range("N136").Value = -55016.073

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 508.69232
$ws.Range("I4").Value = 234.77777
$ws.Range("J4").Value = 1125
$ws.Range("K4").Value = 234.77777
$ws.Range("L4").Value = 1125
$ws.Range("M4").Value = -121.77777
$ws.Range("N4").Value = -1351
$ws.Range("H34").Value = 16000
$ws.Range("I34").Value = 12000
$ws.Range("K34").Value = 12000
$ws.Range("M34").Value = -11797
$ws.Range("H62").Value = 7504.647
$ws.Range("J62").Value = 5582.7144
$ws.Range("L62").Value = 5582.7144
$ws.Range("N62").Value = -6830.7144
$ws.Range("H65").Value = 7504.647
$ws.Range("J65").Value = 5582.7144
$ws.Range("L65").Value = 27913.572
$ws.Range("N65").Value = -34153.572
$ws.Range("H104").Value = 26356.334
$ws.Range("J104").Value = 26356.334
$ws.Range("L104").Value = 26356.334
$ws.Range("N104").Value = -33344.334
